# "cleaned up plots with viridis" - update CCB project management tracker:
#  - rename project_subtype "Differential Expression + GSEA" -> "DE + GSEA"
#  - rename project_type "CRISPR-Cas9" -> "CRISPR screen"
#  - row 22 (MAGeCK-VISPR CRISPR-Cas Analysis): project_subtype CRISPR-Cas9 -> Epigenetics
#  - duration_hrs updates on a few rows
#  - cosmetic: narrower column C, selection moved to B22

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename "Differential Expression + GSEA" -> "DE + GSEA" wherever it appears in column C ---
$deGseaRows = @(2, 8, 9, 18, 21)
foreach ($r in $deGseaRows) {
    $ws.Cells.Item($r, 3).Value = "DE + GSEA"
}

# --- Rename "CRISPR-Cas9" -> "CRISPR screen" wherever it appears in column B (project_type) ---
$crisprTypeRows = @(6, 16, 22)
foreach ($r in $crisprTypeRows) {
    $ws.Cells.Item($r, 2).Value = "CRISPR screen"
}

# --- Row 22: project_subtype changes from CRISPR-Cas9 to Epigenetics ---
$ws.Cells.Item(22, 3).Value = "Epigenetics"

# --- duration_hrs (column G) corrections ---
$ws.Cells.Item(2, 7).Value = 140
$ws.Cells.Item(19, 7).Value = 16
$ws.Cells.Item(22, 7).Value = 25
$ws.Cells.Item(27, 7).Value = 3

# --- cosmetic changes ---
# Note: Excel's ColumnWidth setter round-trips through a pixel-based MDW
# conversion that adds ~5/6 of a character unit when written back out as
# the raw OOXML <col width>. Compensate so the saved width lands on 17.5.
$ws.Columns.Item(3).ColumnWidth = (17.5 - 5/6)
$ws.Range("B22").Select()
